# sp2 review and updates
#
# 1) The "date last updated" field cached on the Slide Master, every
#    Custom Layout, and the Notes Master moves on from 16/08/2018 to
#    21/08/2018.
# 2) Slide 1's agenda bullet for "Slide 8" is reworded.

$p = $ppt.ActivePresentation

$oldDate = "16/08/2018"
$newDate = "21/08/2018"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide Master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every Custom (slide) Layout hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Notes Master
Update-DatePlaceholder $p.NotesMaster.Shapes

# Slide 1: reword the "Slide 8" agenda bullet
$slide1 = $p.Slides.Item(1)
$contentShape = $slide1.Shapes.Item("Content Placeholder 2")
$tr = $contentShape.TextFrame.TextRange
$oldBullet = " " + [char]0x2013 + " Review of the IA report (this might not be finished yet)."
$newBullet = " " + [char]0x2013 + " Review of the your card sort and tree jack exercises"
$found = $tr.Find($oldBullet, 0)
if ($found) {
    $found.Text = $newBullet
}
